# Fixed naive component forecaster bug - Presentation state 11.02.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New matched-error values per row (columns B..K), $null = leave cell blank
$data = @{
    2 = @(1.056122842320784, 0.4458392755916881, 0.35701661626763, 0.08641341155685073, 0.6769318707313746, 0.3469208787210216, 0.3687093173527261, 0.4762867280777406, 0.6850938218533179, 0.1817482544735015)
    3 = @(0.4682672552246967, 0.3596201558300787, 0.06434572590283549, 0.6709498904270196, 0.3415776285386105, 0.3583407839823598, 0.4678094059541168, 0.6773215135465072, 0.1731242809405329, 0.4653169537672726)
    4 = @(0.3810447272872988, 0.1845235747020479, 0.5817555460132365, 0.3118506246900327, 0.3669669040122699, 0.4418692325286816, 0.6545429505182796, 0.1595062457919151, 0.4467961565820968, 0.1882892220727796)
    5 = @(0.1422029423816584, 0.5525510721575033, 0.3088300248851115, 0.3481056204007895, 0.4215381014634773, 0.6400127950840317, 0.1428652002996471, 0.4292123111533466, 0.1718153385427368, 0.5550469433309027)
    6 = @(0.8927488127754134, 0.3842010866690486, 0.1570123340877904, 0.4462320140502105, 0.6493921986403277, 0.08805440168825607, 0.4103182202646768, 0.1569513377276911, 0.526335232701532, 0.1198157454308724)
    7 = @(0.8348959586592991, 0.2041230363001488, 0.2059487832594023, 0.6853372797061905, 0.08600364424866319, 0.3337499808286466, 0.1295519433524077, 0.5000461825038066, 0.0764184547175398, 0.3733815035592049)
    8 = @(0.5164486232236872, 0.3392375229949266, 0.506258857889999, 0.1141978830192304, 0.3697752920210401, 0.09201037314819407, 0.4945038434164454, 0.08059378365628841, 0.3621967898194339, $null)
    9 = @(0.5748280141027678, 0.5909375987643086, -0.03140248361448672, 0.3797290348802828, 0.107375542847739, 0.4556555084590223, 0.06474797462910251, 0.3538625613519331, $null, $null)
    10 = @(0.9019617852456914, 0.08568629079670848, 0.2172648934307159, 0.1363344933129406, 0.4913983809139259, 0.03342032824547542, 0.3506375116504074, $null, $null, $null)
    11 = @(0.3325114682008229, 0.2349207609686054, 0.0417122542019461, 0.5236157691624059, 0.0458084540731927, 0.3297737081011285, $null, $null, $null, $null)
    12 = @(0.4745971203848173, 0.1267364976711596, 0.4071311908043919, 0.06091372572504519, 0.3530127437271186, $null, $null, $null, $null, $null)
    13 = @(0.2915965747052469, 0.4206684630523081, -0.002059193264917797, 0.3630509094086165, $null, $null, $null, $null, $null, $null)
    14 = @(0.6745175049177161, 0.09725885691711864, 0.2514743671933744, $null, $null, $null, $null, $null, $null, $null)
    15 = @(0.1413185481565676, 0.2728948720679752, $null, $null, $null, $null, $null, $null, $null, $null)
    16 = @(0.5091380033804217, $null, $null, $null, $null, $null, $null, $null, $null, $null)
}

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    for ($i = 0; $i -lt $rowVals.Length; $i++) {
        $col = $i + 2   # column B = 2
        $val = $rowVals[$i]
        if ($null -ne $val) {
            $ws.Cells.Item($r, $col).Value = $val
        }
    }
}

